$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 ("HighKick") gains "ready to be fit" markers in the same columns
# that row 11/9 already use (U/V), plus several more (C, H, I, J, K, S, T).
$cells = @("C12", "H12", "I12", "J12", "K12", "S12", "T12", "U12", "V12")
foreach ($cellRef in $cells) {
    $ws.Range($cellRef).Value = "ready to be fit"
}

# Column C widened from 29 to 38.5 characters.
$ws.Columns("C").ColumnWidth = 37.625

# Active cell/selection moved from C31 to C17.
$ws.Range("C17").Select() | Out-Null
